$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the content/description text for the 8 joulu entry (row 43)
$ws.Range("C43").Value = "Fluiditutoriaalin palastelua, rajaava voluumi ja alkutoimia"

# Update the time range text for the same entry
$ws.Range("B43").Value = "17.45-19.45, 20.00-21.00"

# Add the hours logged for this entry, which feeds into the H3 = SUM(G3:G60) total
$ws.Range("G43").Value = 3

# Row 43 grows taller to fit the now-longer wrapped text
$ws.Rows.Item(43).RowHeight = 28.8

# Move the active selection to H43
$ws.Range("H43").Select()

$wb.Save()
